$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.224.82"

$ws.Range("D3").Value = "1.909.12"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  -0.34%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "326.46"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.47%  "

$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("E7").Value = "  +0.36%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.4029"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +2.38%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.08474"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "42.70"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.117"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.11%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "23.43"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +13.44%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "6.461"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.27%  "

$ws.Range("D14").Value = "1.927.87"
$ws.Range("E14").Value = "  +1.23%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.365"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("E16").Value = "  -0.40%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "95.17"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +2.13%  "

$ws.Range("E18").Value = "  +0.76%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06686"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "18.38"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +2.56%  "

$ws.Range("E21").Value = "  -0.33%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.998"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.54%  "

$ws.Range("D23").Value = "30.222.22"
$ws.Range("E23").Value = "  +3.24%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.29"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.33%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.212"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").Value = "2.126.84"
$ws.Range("E26").Value = "  +0.28%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "21.76"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +4.01%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "161.68"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

$ws.Range("E29").Value = "  -2.07%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "129.60"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +1.88%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.098"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +3.66%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.1060"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +1.12%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "6.069"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.759"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +3.11%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.02506"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.09%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.06588"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.2220"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.12%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "5.259"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +2.66%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.232"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "11.92"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +5.72%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "8.814"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.53%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.6526"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.16%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.234"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.6130"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "13.28"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.720"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.16%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.065"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.246"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.28%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "125.01"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.64%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.160"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "79.53"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +2.12%  "
